$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# The model sheet used to carry a generic JSON-ish type (string/number/object)
# plus a separate elementType column (C) for nested/complex prompts (e.g.
# geopoint stored as an "object"). Now the "type" column should just hold the
# actual prompt type (matching the survey sheet's type column directly,
# recursively), so the elementType column is no longer needed.

# Update the "type" values (column B) to the real prompt types.
$ws.Range("B2").Value = "text"        # refrigerator_id
$ws.Range("B3").Value = "decimal"     # refrigerator_size
$ws.Range("B4").Value = "geopoint"    # refrigerator_location (was object/elementType=geopoint)
$ws.Range("B5").Value = "select_one"  # refrigerator_condition

# Drop the now-unneeded elementType column (C).
$ws.Columns.Item(3).Delete()

# Move the active selection to B4, matching where the edit was made.
[void]$ws.Range("B4").Select()
